$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Kingdom"
$ws.Range("E1").Value = "Phylum"
$ws.Range("F1").Value = "Class"
$ws.Range("G1").Value = "Order"
$ws.Range("H1").Value = "Family"
$ws.Range("I1").Value = "Genus"
$ws.Range("J1").Value = "species"

$ws.Range("C1").Copy()
$ws.Range("D1:J1").PasteSpecial(-4122)

$ws.Range("D4").Value = "Unassigned"
$ws.Range("E4").Value = "Unassigned"
$ws.Range("F4").Value = "Unassigned"
$ws.Range("G4").Value = "Unassigned"
$ws.Range("H4").Value = "Unassigned"
$ws.Range("I4").Value = "Unassigned"
$ws.Range("J4").Value = "Unassigned"
